# Update betting odds values for the two match rows (row 2: Guarani vs Amazonas,
# row 3: Botafogo SP vs Ceara) to reflect the latest FlashScore snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("G2").Value = 1.95
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 4.1
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 2.38
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("AC2").Value = 7
$ws.Range("AJ2").Value = 41
$ws.Range("AP2").Value = 26
$ws.Range("AT2").Value = 2.38

# --- Row 3 ---
$ws.Range("G3").Value = 4.5
$ws.Range("H3").Value = 3.1
$ws.Range("K3").Value = 1.95
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("W3").Value = 9.5
$ws.Range("X3").Value = 21
$ws.Range("Z3").Value = 51
$ws.Range("AC3").Value = 6.5
$ws.Range("AF3").Value = 67
$ws.Range("AK3").Value = 19
$ws.Range("AX3").Value = 11
$ws.Range("AY3").Value = 26
$ws.Range("BA3").Value = 67

$wb.Save()
